# Update the "取得日時" (retrieved datetime) column for all existing data rows
# on the "ランサーズ" sheet to reflect the new run timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-03 18:32:56"

# Data rows 2 through 21 all share the same retrieval timestamp; bump them
# to the latest run time while keeping the values as plain text.
for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
